$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to be written as literal text (avoids Excel auto-converting
# numeric-looking strings such as "300.07" into floating point numbers).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
$ws.Range("D2").Value = "43.006.95"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "2.305.08"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "300.07"

# Row 6
Set-TextValue $ws.Range("D6") "97.88"
$ws.Range("E6").Value = "  -0.75%  "

# Row 7
$ws.Range("E7").Value = "  -1.91%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.507"
$ws.Range("E9").Value = "  -2.53%  "

# Row 10
$ws.Range("E10").Value = "  +0.80%  "

# Row 11
$ws.Range("E11").Value = "  -0.04%  "

# Row 12
Set-TextValue $ws.Range("D12") "18.16"
$ws.Range("E12").Value = "  +1.60%  "

# Row 13
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$ws.Range("D15").Value = "2.664.57"
$ws.Range("E15").Value = "  +0.04%  "

# Row 16
$ws.Range("D16").Value = "2.301.48"
$ws.Range("E16").Value = "  -1.88%  "

# Row 17
$ws.Range("E17").Value = "  -0.88%  "

# Row 18
$ws.Range("D18").Value = "42.932.25"
$ws.Range("E18").Value = "  -0.12%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.71"
$ws.Range("E19").Value = "  -5.87%  "

# Row 20
$ws.Range("E20").Value = "  -0.41%  "

# Row 21
$ws.Range("E21").Value = "  -1.75%  "

# Row 22
Set-TextValue $ws.Range("D22") "67.95"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23
Set-TextValue $ws.Range("D23") "235.97"
$ws.Range("E23").Value = "  -1.32%  "

# Row 24
$ws.Range("E24").Value = "  -1.53%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  +0.20%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D26") "2.45"
$ws.Range("E26").Value = "  +0.49%  "

# Row 27
$ws.Range("E27").Value = "  -0.82%  "

# Row 28
Set-TextValue $ws.Range("D28") "25.36"
$ws.Range("E28").Value = "  +2.34%  "

# Row 29
Set-TextValue $ws.Range("D29") "165.72"
$ws.Range("E29").Value = "  -1.23%  "

# Row 30
$ws.Range("E30").Value = "  +0.37%  "

# Row 31
$ws.Range("E31").Value = "  -1.23%  "

# Row 32
Set-TextValue $ws.Range("D32") "33.24"
$ws.Range("E32").Value = "  -0.45%  "

# Row 33
Set-TextValue $ws.Range("D33") "4.90"
$ws.Range("E33").Value = "  +1.47%  "

# Row 34
$ws.Range("E34").Value = "  +0.05%  "

# Row 35
$ws.Range("E35").Value = "  -3.97%  "

# Row 36
Set-TextValue $ws.Range("D36") "17.01"
$ws.Range("E36").Value = "  -6.40%  "

# Row 37
$ws.Range("E37").Value = "  -1.16%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.0687"
$ws.Range("E38").Value = "  -0.53%  "

# Row 39
$ws.Range("E39").Value = "  -0.92%  "

# Row 40
$ws.Range("E40").Value = "  -1.78%  "

# Row 41
$ws.Range("E41").Value = "  -1.81%  "

# Row 42
$ws.Range("E42").Value = "  -0.64%  "

# Row 43
$ws.Range("D43").Value = "2.014.47"
$ws.Range("E43").Value = "  +0.47%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0282"
$ws.Range("E44").Value = "  -1.96%  "

# Row 45
$ws.Range("E45").Value = "  +0.21%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.12"
$ws.Range("E46").Value = "  -1.64%  "

# Row 47
Set-TextValue $ws.Range("D47") "17.65"
$ws.Range("E47").Value = "  +1.24%  "

# Row 48
$ws.Range("E48").Value = "  -1.37%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.98"
$ws.Range("E49").Value = "  +1.49%  "

# Row 50
Set-TextValue $ws.Range("D50") "53.90"
$ws.Range("E50").Value = "  -1.12%  "

# Row 51
$ws.Range("D51").Value = "2.532.12"
$ws.Range("E51").Value = "  +0.12%  "
